$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '44.239.54'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +2.26%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.379.25'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('E4').Value = '  -0.19%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '0.694'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +7.47%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '243.29'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +4.29%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '76.84'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +6.79%  '
$ws.Range('E8').Value = '  -0.12%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.610'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +28.04%  '
$ws.Range('E10').Value = '  +7.23%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '57.96'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +1.98%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '32.60'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +20.40%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '7.59'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +21.12%  '
$ws.Range('E14').Value = '  +2.50%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '2.733.28'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.24%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '17.19'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +8.05%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.930'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +8.84%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.363.97'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.13%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '44.395.02'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('E20').Value = '  +4.78%  '
$ws.Range('E21').Value = '  +6.98%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '78.93'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +6.40%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '257.69'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +3.63%  '
$ws.Range('E24').Value = '  +0.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.59'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +5.42%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '3.71'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '10.96'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +10.34%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.79'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +19.12%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +2.09%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '23.17'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +3.97%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '175.64'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('E33').Value = '  +7.80%  '
$ws.Range('E34').Value = '  +8.68%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0762'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +10.73%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.36'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +6.44%  '
$ws.Range('E37').Value = '  +5.90%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.51'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +2.94%  '
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  +9.70%  '
$ws.Range('E41').Value = '  +3.27%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '19.09'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +3.79%  '
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  +18.60%  '
$ws.Range('E45').Value = '  +3.93%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +13.99%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.28'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +5.90%  '
$ws.Range('E48').Value = '  +6.31%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '103.18'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +3.78%  '
$ws.Range('E50').Value = '  -0.83%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.479.00'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +2.66%  '
